$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster order for column B (rows 2-32), with $null meaning "clear the cell"
$values = @(
    "Yunjae",
    "遠藤隼人",
    "富澤天音",
    "神山修造",
    "志塚惇希",
    "富澤天音",
    "豊島亮",
    $null,
    "兒島大志郎",
    "日高泰聖",
    "白岩詩佑介",
    "Cox Matthew Jonah",
    "Hansen Jakob U",
    "石井海成",
    "Nicholas Tristan Aryasatyo",
    "小溝賢",
    "小野文哉",
    "渡部魁",
    "崎谷航平",
    "三神佳誠",
    "氏家琉貴",
    "羽賀尚生",
    "島田実",
    "足立耕平",
    "Yunjae",
    "富澤天音",
    "神山修造",
    "志塚惇希",
    $null,
    $null,
    $null
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $val = $values[$i]
    if ($null -eq $val) {
        $ws.Cells.Item($row, 2).Value = ""
    } else {
        $ws.Cells.Item($row, 2).Value = $val
    }
}

$ws.Range("D7").Select()
